# Auto-generated edit script updating the cryptos price list: refreshed
# Price/Volume(1h) figures and a swap of the Stellar/BinanceUSD rows (28 & 29),
# matching the upstream GitHub Actions data-refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.576.00'
$ws.Range("E2").Value = '  -1.37%  '
$ws.Range("D3").Value = '1.666.73'
$ws.Range("E3").Value = '  -3.40%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.17'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.65%  '
$ws.Range("E6").Value = '  -1.76%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.55'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.98%  '
$ws.Range("E9").Value = '  -1.30%  '
$ws.Range("E10").Value = '  -1.90%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0880'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.33%  '
$ws.Range("D12").Value = '1.902.39'
$ws.Range("D13").Value = '1.659.02'
$ws.Range("E13").Value = '  -3.70%  '
$ws.Range("E14").Value = '  -2.61%  '
$ws.Range("E15").Value = '  -2.18%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.29'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '250.92'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.63%  '
$ws.Range("D18").Value = '27.592.62'
$ws.Range("E18").Value = '  -1.12%  '
$ws.Range("D19").Value = '0.0₃0732'
$ws.Range("E19").Value = '  -3.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.54'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -4.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.999'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.00%  '
$ws.Range("E22").Value = '  -3.02%  '
$ws.Range("E23").Value = '  -4.59%  '
$ws.Range("E24").Value = '  -5.61%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.50'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.87%  '
$ws.Range("E26").Value = '  -1.47%  '
$ws.Range("E27").Value = '  -5.03%  '
$ws.Range("B28").Value = 'BinanceUSD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.112'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.34%  '
$ws.Range("E30").Value = '  +4.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0508'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.78%  '
$ws.Range("E32").Value = '  -2.65%  '
$ws.Range("D33").Value = '1.475.91'
$ws.Range("E33").Value = '  -1.13%  '
$ws.Range("E34").Value = '  -5.49%  '
$ws.Range("E35").Value = '  -5.32%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.943'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.75%  '
$ws.Range("E37").Value = '  -0.94%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.577'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -5.95%  '
$ws.Range("E39").Value = '  -2.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '69.77'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.48%  '
$ws.Range("E41").Value = '  -3.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("E43").Value = '  -6.98%  '
$ws.Range("D44").Value = '1.810.73'
$ws.Range("E44").Value = '  -3.33%  '
$ws.Range("E45").Value = '  -3.25%  '
$ws.Range("E46").Value = '  -0.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.70'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '89.52'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.71%  '
$ws.Range("E49").Value = '  -2.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '42.21'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +15.97%  '
$ws.Range("E51").Value = '  -3.56%  '
